# Applies the "Elimina EC anteriores y se agregan nuevos, se modifica base de datos" edit:
#  - Reorders / refreshes the debtor periods for LUIS MIGUEL CORTECERO HURTADO
#  - Replaces the second worker's data (was ALY ANDREA BARRIOS CANTILLO 2211/2212)
#    with DANIEL ENRIQUE GUERRERO VALENZUELA (period 2202)
#  - Re-adds ALY ANDREA BARRIOS CANTILLO (periods 2212, 2211) further down the table
#  - Adds a new worker JORGE LUIS PAJARO VALENZUELA (periods 2507, 2506)
#  - Refreshes the summary totals (VALOR MORA, Cant. Trabajadores, Cant. Periodos)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Make room for the two extra detail rows needed (table grows from 9 to 12 rows) ---
# Row 24 currently holds the last (closing-style) detail row; insert 3 blank rows above it
# so it ends up at row 27, and copy the "normal" row formatting into the new rows 24-26.
$ws.Rows("24:26").Insert()
$ws.Range("B23:J23").Copy() | Out-Null
$ws.Range("B24:J26").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Worker 1: LUIS MIGUEL CORTECERO HURTADO (CC 1044917994), periods refreshed ---
$ws.Range("E16").Value = "2110"
$ws.Range("F16").Value = 52000
$ws.Range("G16").Value = 1300000

$ws.Range("E17").Value = "2109"
$ws.Range("F17").Value = 52000
$ws.Range("G17").Value = 1300000

$ws.Range("E18").Value = "2108"
$ws.Range("F18").Value = 52000
$ws.Range("G18").Value = 1300000

$ws.Range("E19").Value = "2107"
$ws.Range("F19").Value = 52000
$ws.Range("G19").Value = 1300000

$ws.Range("E20").Value = "2106"
$ws.Range("F20").Value = 52000
$ws.Range("G20").Value = 1300000

$ws.Range("E21").Value = "2105"
$ws.Range("F21").Value = 52000
$ws.Range("G21").Value = 1300000

$ws.Range("E22").Value = "2104"
$ws.Range("F22").Value = 43333
$ws.Range("G22").Value = 1300000

# --- Worker 2: DANIEL ENRIQUE GUERRERO VALENZUELA (CC 1051444016), period 2202 ---
$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = "1051444016"
$ws.Range("D23").Value = "DANIEL ENRIQUE GUERRERO VALENZUELA"
$ws.Range("E23").Value = "2202"
$ws.Range("F23").Value = 22400
$ws.Range("G23").Value = 1000000

# --- Worker 3: ALY ANDREA BARRIOS CANTILLO (CC 1235039727), periods 2212 & 2211 ---
$ws.Range("B24").Value = "CC"
$ws.Range("C24").Value = "1235039727"
$ws.Range("D24").Value = "ALY ANDREA BARRIOS CANTILLO"
$ws.Range("E24").Value = "2212"
$ws.Range("F24").Value = 64000
$ws.Range("G24").Value = 1600000

$ws.Range("B25").Value = "CC"
$ws.Range("C25").Value = "1235039727"
$ws.Range("D25").Value = "ALY ANDREA BARRIOS CANTILLO"
$ws.Range("E25").Value = "2211"
$ws.Range("F25").Value = 21333
$ws.Range("G25").Value = 1600000

# --- Worker 4: JORGE LUIS PAJARO VALENZUELA (CC 1051447433), periods 2507 & 2506 ---
$ws.Range("B26").Value = "CC"
$ws.Range("C26").Value = "1051447433"
$ws.Range("D26").Value = "JORGE LUIS PAJARO VALENZUELA"
$ws.Range("E26").Value = "2507"
$ws.Range("F26").Value = 56940
$ws.Range("G26").Value = 1423500

$ws.Range("B27").Value = "CC"
$ws.Range("C27").Value = "1051447433"
$ws.Range("D27").Value = "JORGE LUIS PAJARO VALENZUELA"
$ws.Range("E27").Value = "2506"
$ws.Range("F27").Value = 56940
$ws.Range("G27").Value = 1423500

# --- Refresh summary block ---
$ws.Range("E11").Value = 576946
$ws.Range("C13").Value = 4
$ws.Range("F13").Value = 12

# --- Column D needs to be a bit wider to fit the longer new names ---
$ws.Columns("D").AutoFit()
